$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new daily records (row 18: 四方坪站, row 19: 高岭站) for date 45939,
# mirroring the formatting of the preceding rows.

$ws.Range("A17:F17").Copy()
$ws.Range("A18:F19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(18, 1).Value = 45939
$ws.Cells.Item(18, 2).Value = "四方坪站"
$ws.Cells.Item(18, 3).Value = 10987.23
$ws.Cells.Item(18, 4).Value = 9022.8700000000008
$ws.Cells.Item(18, 5).Value = 3818.52
$ws.Cells.Item(18, 6).Value = 460

$ws.Cells.Item(19, 1).Value = 45939
$ws.Cells.Item(19, 2).Value = "高岭站"
$ws.Cells.Item(19, 3).Value = 5141.28
$ws.Cells.Item(19, 4).Value = 4172.6499999999996
$ws.Cells.Item(19, 5).Value = 1323.75
$ws.Cells.Item(19, 6).Value = 179

$ws.Range("E21").Select()
